$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename worksheet to English name
$ws.Name = "Orders"

# Translate header row (row 1) from Arabic to English
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "Customer Name"
$ws.Range("C1").Value = "Status"
$ws.Range("D1").Value = "Created At"
$ws.Range("E1").Value = "Updated At"
$ws.Range("F1").Value = "Phone Number"
$ws.Range("G1").Value = "Street Name"
$ws.Range("H1").Value = "Province"
$ws.Range("J1").Value = "Shipping Cost"
$ws.Range("K1").Value = "Notes"
$ws.Range("L1").Value = "Country"
$ws.Range("M1").Value = "VAT Profit"
$ws.Range("N1").Value = "Order Profit"
$ws.Range("O1").Value = "Page Name"
$ws.Range("P1").Value = "Page URL"
$ws.Range("Q1").Value = "SKUs"
$ws.Range("R1").Value = "Quantities"
$ws.Range("S1").Value = "Prices"
$ws.Range("T1").Value = "Order Received By"
$ws.Range("U1").Value = "National Address"
$ws.Range("V1").Value = "Source of the National Address"
$ws.Range("W1").Value = "Store Order ID"

# Translate the "Status" column (C) for every data row (2-147)
$ws.Range("C2").Value = "Order Received"
$ws.Range("C3").Value = "Order Confirmed"
$ws.Range("C4").Value = "Customer Rejected"
$ws.Range("C5").Value = "Order Received"
$ws.Range("C6").Value = "Order Confirmed"
$ws.Range("C7").Value = "Customer Rejected"
$ws.Range("C8").Value = "Order Confirmed"
$ws.Range("C9").Value = "Order Received"
$ws.Range("C10").Value = "Delivery in Progress"
$ws.Range("C11").Value = "Delivery in Progress"
$ws.Range("C12").Value = "Delivery in Progress"
$ws.Range("C13").Value = "Pending Shipping"
$ws.Range("C14").Value = "Delivery in Progress"
$ws.Range("C15").Value = "Pending Shipping"
$ws.Range("C16").Value = "Pending Shipping"
$ws.Range("C17").Value = "Delivery in Progress"
$ws.Range("C18").Value = "Delivery in Progress"
$ws.Range("C19").Value = "Order Received"
$ws.Range("C20").Value = "Customer Rejected"
$ws.Range("C21").Value = "Customer Rejected"
$ws.Range("C22").Value = "Delivery in Progress"
$ws.Range("C23").Value = "Delivered"
$ws.Range("C24").Value = "Delivered"
$ws.Range("C25").Value = "Order Received"
$ws.Range("C26").Value = "Customer Rejected"
$ws.Range("C27").Value = "Delivered"
$ws.Range("C28").Value = "Delivered"
$ws.Range("C29").Value = "Delivery in Progress"
$ws.Range("C30").Value = "Delivery in Progress"
$ws.Range("C31").Value = "Customer Rejected"
$ws.Range("C32").Value = "Delivery in Progress"
$ws.Range("C33").Value = "Customer Rejected"
$ws.Range("C34").Value = "Order Confirmed"
$ws.Range("C35").Value = "Cancelled by You"
$ws.Range("C36").Value = "Delivered"
$ws.Range("C37").Value = "Cancelled by You"
$ws.Range("C38").Value = "Delivery in Progress"
$ws.Range("C39").Value = "Delivery in Progress"
$ws.Range("C40").Value = "Delivered"
$ws.Range("C41").Value = "Delivery in Progress"
$ws.Range("C42").Value = "Delivered"
$ws.Range("C43").Value = "Delivered"
$ws.Range("C44").Value = "Customer Rejected"
$ws.Range("C45").Value = "Delivered"
$ws.Range("C46").Value = "Delivery in Progress"
$ws.Range("C47").Value = "Delivered"
$ws.Range("C48").Value = "Delivery in Progress"
$ws.Range("C49").Value = "Delivery in Progress"
$ws.Range("C50").Value = "Delivered"
$ws.Range("C51").Value = "Order Received"
$ws.Range("C52").Value = "Delivered"
$ws.Range("C53").Value = "Cancelled by You"
$ws.Range("C54").Value = "Customer Rejected"
$ws.Range("C55").Value = "Delivered"
$ws.Range("C56").Value = "Delivered"
$ws.Range("C57").Value = "Cancelled by You"
$ws.Range("C58").Value = "Delivered"
$ws.Range("C59").Value = "Delivery Failed"
$ws.Range("C60").Value = "Delivery Failed"
$ws.Range("C61").Value = "Delivered"
$ws.Range("C62").Value = "Customer Rejected"
$ws.Range("C63").Value = "Cancelled by You"
$ws.Range("C64").Value = "Cancelled by You"
$ws.Range("C65").Value = "Cancelled by You"
$ws.Range("C66").Value = "Customer Rejected"
$ws.Range("C67").Value = "Cancelled by You"
$ws.Range("C68").Value = "Customer Rejected"
$ws.Range("C69").Value = "Delivery in Progress"
$ws.Range("C70").Value = "Delivery in Progress"
$ws.Range("C71").Value = "Delivered"
$ws.Range("C72").Value = "Cancelled by You"
$ws.Range("C73").Value = "Delivery in Progress"
$ws.Range("C74").Value = "Cancelled by You"
$ws.Range("C75").Value = "Delivered"
$ws.Range("C76").Value = "Delivered"
$ws.Range("C77").Value = "Delivered"
$ws.Range("C78").Value = "Customer Rejected"
$ws.Range("C79").Value = "Customer Rejected"
$ws.Range("C80").Value = "Delivered"
$ws.Range("C81").Value = "Customer Rejected"
$ws.Range("C82").Value = "Customer Rejected"
$ws.Range("C83").Value = "Customer Rejected"
$ws.Range("C84").Value = "Delivery in Progress"
$ws.Range("C85").Value = "Delivered"
$ws.Range("C86").Value = "Cancelled by You"
$ws.Range("C87").Value = "Cancelled by You"
$ws.Range("C88").Value = "Cancelled by You"
$ws.Range("C89").Value = "Cancelled by You"
$ws.Range("C90").Value = "Delivered"
$ws.Range("C91").Value = "Delivered"
$ws.Range("C92").Value = "Delivered"
$ws.Range("C93").Value = "Cancelled by You"
$ws.Range("C94").Value = "Cancelled by You"
$ws.Range("C95").Value = "Cancelled by You"
$ws.Range("C96").Value = "Cancelled by You"
$ws.Range("C97").Value = "Customer Rejected"
$ws.Range("C98").Value = "Customer Rejected"
$ws.Range("C99").Value = "Cancelled by You"
$ws.Range("C100").Value = "Cancelled by You"
$ws.Range("C101").Value = "Cancelled by You"
$ws.Range("C102").Value = "Delivered"
$ws.Range("C103").Value = "Delivery in Progress"
$ws.Range("C104").Value = "Delivered"
$ws.Range("C105").Value = "Delivered"
$ws.Range("C106").Value = "Cancelled by You"
$ws.Range("C107").Value = "Cancelled by You"
$ws.Range("C108").Value = "Cancelled by You"
$ws.Range("C109").Value = "Delivered"
$ws.Range("C110").Value = "Delivered"
$ws.Range("C111").Value = "Customer Rejected"
$ws.Range("C112").Value = "Delivered"
$ws.Range("C113").Value = "Customer Rejected"
$ws.Range("C114").Value = "Customer Rejected"
$ws.Range("C115").Value = "Delivered"
$ws.Range("C116").Value = "Customer Rejected"
$ws.Range("C117").Value = "Delivered"
$ws.Range("C118").Value = "Customer Rejected"
$ws.Range("C119").Value = "Customer Rejected"
$ws.Range("C120").Value = "Cancelled by You"
$ws.Range("C121").Value = "Return Verified"
$ws.Range("C122").Value = "Delivered"
$ws.Range("C123").Value = "Delivered"
$ws.Range("C124").Value = "Delivered"
$ws.Range("C125").Value = "Cancelled by You"
$ws.Range("C126").Value = "Cancelled by You"
$ws.Range("C127").Value = "Cancelled by You"
$ws.Range("C128").Value = "Customer Rejected"
$ws.Range("C129").Value = "Cancelled by You"
$ws.Range("C130").Value = "Return Verified"
$ws.Range("C131").Value = "Temporary Suspended"
$ws.Range("C132").Value = "Delivered"
$ws.Range("C133").Value = "Cancelled by You"
$ws.Range("C134").Value = "Delivered"
$ws.Range("C135").Value = "Customer Rejected"
$ws.Range("C136").Value = "Customer Rejected"
$ws.Range("C137").Value = "Return Verified"
$ws.Range("C138").Value = "Customer Rejected"
$ws.Range("C139").Value = "Customer Rejected"
$ws.Range("C140").Value = "Delivered"
$ws.Range("C141").Value = "Customer Rejected"
$ws.Range("C142").Value = "Cancelled by You"
$ws.Range("C143").Value = "Cancelled by You"
$ws.Range("C144").Value = "Return Verified"
$ws.Range("C145").Value = "Delivered"
$ws.Range("C146").Value = "Customer Rejected"
$ws.Range("C147").Value = "Customer Rejected"

# Translate the "Country" column (L) for every data row -> all become "Iraq"
$ws.Range("L2:L147").Value = "Iraq"

# Rows 14 and 18 also received a status refresh, which updated their "Updated At" timestamp
$ws.Range("E14").Value = 46068.927126145834
$ws.Range("E18").Value = 46068.92732148148
